$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: G and H change
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-10: only H changes
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1

# Rows 11-14: D and E change
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1

# Rows 15-16: only H changes
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1

# Row 17: D and E change
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18: only H changes
$ws.Range("H18").Value = 1
